$d = $word.ActiveDocument

# 1. "Add a README file to the upload." -> "Add a README file (in TXT format) to the upload."
$d.Content.Find.Execute(
    "Add a README file to the upload.", $true, $false, $false, $false, $false,
    $true, 1, $false, "Add a README file (in TXT format) to the upload.", 2)

# 2-4. "Specify in the README:" -> "Specify in the README (in TXT format):"
# (3 occurrences in the document) - wdReplaceAll (2) replaces every match in one call.
$d.Content.Find.Execute(
    "Specify in the README:", $true, $false, $false, $false, $false,
    $true, 1, $false, "Specify in the README (in TXT format):", 2)
